# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
#
# Two pairs of adjacent match rows had their data (everything except the
# leading running-index column A) swapped between rows. Column A keeps the
# row's original running index; columns B:AB (match id, teams, score,
# result, and all odds columns) swap between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($sheet, $row1, $row2) {
    $range1 = $sheet.Range("B$row1" + ":AB$row1")
    $range2 = $sheet.Range("B$row2" + ":AB$row2")

    $data1 = $range1.Value2
    $data2 = $range2.Value2

    $range1.Value2 = $data2
    $range2.Value2 = $data1
}

Swap-RowData $ws 29 30
Swap-RowData $ws 87 88
Swap-RowData $ws 99 100
Swap-RowData $ws 192 193
